$d = $word.ActiveDocument

function Set-ParagraphPlainText($para, $text) {
    # Collapse every run in the paragraph into a single run carrying $text.
    # Setting Range.Text directly to the already-existing concatenation is a
    # no-op in some engines, so first blank the paragraph body, then write
    # the final text - this reliably merges the runs into one.
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = ""

    $r2 = $para.Range
    $r2.End = $r2.End - 1
    $r2.Text = $text
}

$titlePara = $null
$authorPara = $null
$abstractPara = $null

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Title" -and $titlePara -eq $null) {
        $titlePara = $p
    } elseif ($styleName -eq "Author" -and $authorPara -eq $null) {
        $authorPara = $p
    } elseif ($styleName -eq "Abstract" -and $abstractPara -eq $null) {
        $abstractPara = $p
    }
}

if ($titlePara -ne $null) {
    Set-ParagraphPlainText $titlePara "Questions: Introduction to complex numbers"
}

if ($authorPara -ne $null) {
    Set-ParagraphPlainText $authorPara "Tom Coleman"
}

if ($abstractPara -ne $null) {
    Set-ParagraphPlainText $abstractPara "A selection of questions for the study guide on introduction to complex numbers."
}
